# Apply the ecRxiv indicator metadata edit:
#  - split the combined "NO_URAQ_001-005" indicatorID row into 5 separate rows
#    (NO_URAQ_001 .. NO_URAQ_005), each keeping the same example/reference text
#  - duplicate the "Ecosystem" row 5x and fill in the Ecosystem value
#    (T7.4 Urban and industrial) for each of the 5 indicators
#  - mark the indicator status as "complete" (was "incomplete")
#  - re-anchor the existing hyperlinks to their new row locations

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")

# ---------------------------------------------------------------------------
# 1) Expand row 2 (indicatorID) into 5 rows: insert 4 blank rows below it and
#    copy the formatting across, then fix up the Value (col B) per row.
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt 4; $i++) {
    $ws.Rows.Item(3).Insert()
    $ws.Range("A2:D2").Copy($ws.Range("A3:D3"))
}

$ws.Range("B2").Value2 = "NO_URAQ_001"
$ws.Range("B3").Value2 = "NO_URAQ_002"
$ws.Range("B4").Value2 = "NO_URAQ_003"
$ws.Range("B5").Value2 = "NO_URAQ_004"
$ws.Range("B6").Value2 = "NO_URAQ_005"

$ws.Range("C3").Value2 = "NO_NDVI_002"
$ws.Range("C4").Value2 = "NO_NDVI_003"
$ws.Range("C5").Value2 = "NO_NDVI_004"
$ws.Range("C6").Value2 = "NO_NDVI_005"

# ---------------------------------------------------------------------------
# 2) The "Ecosystem" row used to be row 9; after the 4 rows inserted above it
#    is now row 13. Duplicate it 4x (rows 14-17) and give every one of the 5
#    rows an Ecosystem value.
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt 4; $i++) {
    $ws.Rows.Item(14).Insert()
    $ws.Range("A13:D13").Copy($ws.Range("A14:D14"))
}

$ws.Range("B13").Value2 = "T7.4 Urban and industrial"
$ws.Range("B14").Value2 = "T7.4 Urban and industrial"
$ws.Range("B15").Value2 = "T7.4 Urban and industrial"
$ws.Range("B16").Value2 = "T7.4 Urban and industrial"
$ws.Range("B17").Value2 = "T7.4 Urban and industrial"

for ($r = 13; $r -le 17; $r++) {
    $ws.Rows.Item($r).RowHeight = 28.8
}

# ---------------------------------------------------------------------------
# 3) status: incomplete -> complete (row 12 shifted to row 20)
# ---------------------------------------------------------------------------
$ws.Range("B20").Value2 = "complete"

# ---------------------------------------------------------------------------
# 4) Re-anchor hyperlinks that moved because of the inserted rows.
#    (ECT doi link: old D6 -> D10 ; url row: old B16/C16 -> B24/C24)
# ---------------------------------------------------------------------------
$ectLink = "https://doi.org/10.3897/oneeco.6.e58218"
$urlLink = "https://github.com/NINAnor/ecRxiv/tree/main/indicators/NO_URAQ_000-004"
$exampleUrlLink = "https://github.com/NINAnor/ecRxiv/tree/main/indicators/'INDICATOR-ID'"

$ws.Range("A1:D30").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("D10"), $ectLink)
$ws.Hyperlinks.Add($ws.Range("B24"), $urlLink)
$ws.Hyperlinks.Add($ws.Range("C24"), $exampleUrlLink)

# ---------------------------------------------------------------------------
# 5) Keep the selection / active cell consistent with the new bottom row.
# ---------------------------------------------------------------------------
$ws.Range("B26").Select()
